# Weekly fruit/vegetable price update:
# a new observation is inserted as row 14, pushing the previously-existing
# rows 14-33 down to 15-34 (the data for those rows is otherwise untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14; Excel shifts rows 14:33 down to 15:34
# and carries the row's number formatting (date style on column D) forward.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with this week's record.
$ws.Cells.Item(14, 1).Value  = 9
$ws.Cells.Item(14, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(14, 3).Value  = "Metropolitana"
$ws.Cells.Item(14, 4).Value  = 44708
$ws.Cells.Item(14, 5).Value  = 13
$ws.Cells.Item(14, 6).Value  = 100112035
$ws.Cells.Item(14, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(14, 8).Value  = "Sin especificar"
$ws.Cells.Item(14, 9).Value  = "Primera"
$ws.Cells.Item(14, 10).Value = 25
$ws.Cells.Item(14, 11).Value = 26000
$ws.Cells.Item(14, 12).Value = 26000
$ws.Cells.Item(14, 13).Value = 26000
$ws.Cells.Item(14, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(14, 15).Value = "Hijuelas"
$ws.Cells.Item(14, 16).Value = 1733
$ws.Cells.Item(14, 17).Value = 15
$ws.Cells.Item(14, 18).Value = "Hortaliza"
